# Insert a new weekly price record as row 355 on the "Coliflor" sheet,
# pushing the existing rows 355:439 down to 356:440.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 355 and below down by one (classic Rows.Insert / xlShiftDown).
$ws.Rows(355).Insert()

# Populate the newly inserted row 355 with the new record.
$ws.Cells.Item(355, 1).Value  = 10
$ws.Cells.Item(355, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(355, 3).Value  = "La Araucanía"
$ws.Cells.Item(355, 4).Value  = 44785
$ws.Cells.Item(355, 5).Value  = 9
$ws.Cells.Item(355, 6).Value  = 100112008
$ws.Cells.Item(355, 7).Value  = "Coliflor"
$ws.Cells.Item(355, 8).Value  = "Sin especificar"
$ws.Cells.Item(355, 9).Value  = "Primera"
$ws.Cells.Item(355, 10).Value = 200
$ws.Cells.Item(355, 11).Value = 1500
$ws.Cells.Item(355, 12).Value = 1500
$ws.Cells.Item(355, 13).Value = 1500
$ws.Cells.Item(355, 14).Value = "$/unidad"
$ws.Cells.Item(355, 15).Value = "Región Metropolitana"
$ws.Cells.Item(355, 16).Value = 1500
$ws.Cells.Item(355, 17).Value = 1
$ws.Cells.Item(355, 18).Value = "Hortaliza"
